$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.763.51"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "3.487.21"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "592.55"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "171.47"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.71%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.588"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "4.093.64"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("E13").Value = "  -0.47%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "28.91"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").Value = "66.792.62"
$ws.Range("E15").Value = "  +0.32%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0000177"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "3.499.75"
$ws.Range("E17").Value = "  +0.31%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.26"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.99%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "14.00"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "392.96"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.26%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.94"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "72.96"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  +0.94%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.180"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.84%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.35%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "6.15"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.95%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.42"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("E31").Value = "  -0.35%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "23.67"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.35%  "
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("E34").Value = "  -0.22%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "162.63"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.878"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("E38").Value = "  +3.42%  "
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("E40").Value = "  -0.48%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "27.15"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").Value = "2.816.14"
$ws.Range("E42").Value = "  +1.24%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "26.13"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.33%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "42.80"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.17%  "
$ws.Range("E46").Value = "  -3.60%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "336.28"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.53%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "34.51"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("E49").Value = "  -2.23%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.99%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "6.41"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.23%  "
